$d = $word.ActiveDocument

# 1) Shorten the title paragraph text (keep bold/size formatting, just trim the text).
$d.Content.Find.Execute("A5-1: Critical Thinking Exercises", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A5-1: ", 2) | Out-Null

# 2) Replace the "Complete 3 challenges" paragraph's content + formatting in one shot via raw
#    OOXML (bold/28pt -> bCs/24pt, new body text) using InsertXML so the exact target markup
#    (<w:bCs/>, no <w:b/>) is produced.
$p2 = $d.Paragraphs(2)
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Create an ASP.NET Core (Model-View-Controller) project and take a little tour of all of the files that are created. Try running the project out of the box, and see if you can get the webpage to run!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($newParaXml) | Out-Null

# 3) Delete paragraphs 3 through 10 (1-based: "Calculate average" list item through the
#    "Is this my tail?" image paragraph), keeping the trailing blank paragraph before the
#    section break intact.
$startPara = $d.Paragraphs(3)
$endPara = $d.Paragraphs(10)
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete() | Out-Null
